$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-15-2012-13"
$newDate = "2013-05-15"

# Column BF ("Date") is the 58th column.
$col = 58
$lastRow = $ws.Range("BF1").End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq $oldDate) {
        # Assigning the literal string "2013-05-15" directly via .Value would
        # get auto-recognized as a date and silently coerced into a date
        # serial number. Go through a text formula + paste-as-values round
        # trip instead so the cell keeps holding the literal text
        # "2013-05-15" (matching the source data, which only ever had this
        # value as a plain inline/shared string), without touching the
        # cell's number format/style.
        $cell.Formula = '="' + $newDate + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = $false
